# Swap the full content of paired observation rows (11<->12, 20<->21, 37<->38).
# Each pair of rows in the source sheet had its data reordered; every cell
# from A through AY is exchanged between the two rows in a pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51   # column AY

function Swap-Rows($rowA, $rowB) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        # Value2 reliably reflects the real stored value (numbers as Double,
        # text as String, booleans as Boolean, blanks as empty String) -
        # unlike Value, which on this host returns a bogus reflection token
        # for empty cells.
        $valA = $cellA.Value2
        $valB = $cellB.Value2

        # Only touch cells whose value actually needs to move - this avoids
        # round-tripping text that is identical between the two rows (e.g.
        # date-like text such as "2026-02-20") through the Value setter,
        # which would otherwise get reinterpreted/reformatted even though
        # nothing really changed.
        if ("$valA" -ne "$valB") {
            $cellA.Value = $valB
            $cellB.Value = $valA
        }
    }
}

Swap-Rows 11 12
Swap-Rows 20 21
Swap-Rows 37 38
